$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row => (new Price (column D) or $null, new Volume(1h) (column E))
$updates = @(
    @(2,  "68.183.57", "  +1.16%  "),
    @(3,  "3.564.29",  "  +1.91%  "),
    @(4,  $null,       "  -0.07%  "),
    @(5,  "618.32",    "  +2.57%  "),
    @(6,  "154.51",    "  +3.44%  "),
    @(7,  "3.562.72",  "  +1.90%  "),
    @(8,  $null,       "  -0.04%  "),
    @(9,  $null,       "  +2.26%  "),
    @(10, $null,       "  +5.35%  "),
    @(11, "7.42",      "  +7.05%  "),
    @(12, "0.437",     "  +3.96%  "),
    @(13, "33.22",     "  +5.46%  "),
    @(14, $null,       "  +1.65%  "),
    @(15, "4.167.30",  "  +1.91%  "),
    @(16, "3.561.72",  "  +1.71%  "),
    @(17, "68.247.52", "  +1.40%  "),
    @(18, $null,       "  -0.05%  "),
    @(19, $null,       "  +5.32%  "),
    @(20, $null,       "  +6.71%  "),
    @(21, $null,       "  +11.86%  "),
    @(22, "453.84",    "  +1.48%  "),
    @(23, $null,       "  +4.33%  "),
    @(24, "78.47",     "  +1.44%  "),
    @(25, $null,       "  +1.96%  "),
    @(26, "3.707.79",  "  +1.91%  "),
    @(27, $null,       "  -0.12%  "),
    @(28, "9.23",      "  +12.69%  "),
    @(29, $null,       "  +4.05%  "),
    @(30, $null,       "  +11.47%  "),
    @(31, $null,       "  +3.35%  "),
    @(32, "0.169",     "  +3.31%  "),
    @(33, $null,       "  -0.08%  "),
    @(34, "6.37",      "  +4.94%  "),
    @(35, "26.11",     "  +1.67%  "),
    @(36, $null,       "  +4.90%  "),
    @(37, "3.558.09",  "  +2.00%  "),
    @(38, "8.23",      "  +3.32%  "),
    @(39, $null,       "  +8.70%  "),
    @(40, $null,       "  +0.03%  "),
    @(41, "181.20",    "  +3.56%  "),
    @(42, "0.0916",    "  +4.88%  "),
    @(43, "0.999",     "  -0.05%  "),
    @(44, "5.61",      "  +4.18%  "),
    @(45, "30.88",     "  +12.61%  "),
    @(46, "0.898",     "  +2.16%  "),
    @(47, "46.47",     "  +2.29%  "),
    @(48, "1.33",      "  +5.11%  "),
    @(49, "2.64",      "  +3.66%  "),
    @(50, $null,       "  +3.60%  "),
    @(51, $null,       "  +7.82%  ")
)

foreach ($u in $updates) {
    $row = $u[0]
    $newD = $u[1]
    $newE = $u[2]
    if ($null -ne $newD) {
        $ws.Cells.Item($row, 4).Value = $newD
    }
    $ws.Cells.Item($row, 5).Value = $newE
}
